## Database is added. Database Code is added in PowerPoint.
## Adds a new (3rd) slide containing the four CREATE TABLE statements
## for the database schema (Task, Project, Resource, AssignedTask)
## as separate textboxes, matching the authored slide3.xml.

$p = $ppt.ActivePresentation

# EMU-per-point constant: the Shapes/.Add* family (like real PowerPoint
# automation) takes Left/Top/Width/Height in points, while the target
# OOXML stores EMU (1 pt = 12700 EMU).
$EMU = 12700

# Append a brand-new slide after the existing two, using the same
# "Title and Content" layout as slide 2 (the sibling text-only slide).
$layout = $p.SlideMaster.CustomLayouts.Item(2)
$s = $p.Slides.AddSlide($p.Slides.Count + 1, $layout)

# The layout's Title/Content placeholders aren't used on this slide
# (same as slide 2 in this deck: plain textboxes only); drop them. This
# also naturally advances the shape-id counter so the first real
# textbox below lands on id 4, matching the sibling slides.
while ($s.Shapes.Count -gt 0) {
    $s.Shapes.Item(1).Delete()
}

# --- TextBox 3: Task table (with FK to Project) ------------------------
$tb1 = $s.Shapes.AddTextbox(1, 611560 / $EMU, 2276872 / $EMU, 7704856 / $EMU, 1200329 / $EMU)
$tb1.Name = "TextBox 3"
$tb1.Fill.Visible = 0
$tr1 = $tb1.TextFrame.TextRange
$tr1.Text = "CREATE TABLE Task(TID INTEGER NOT NULL PRIMARY KEY AUTOINCREMENT, PID INT NOT NULL, "
[void]$tr1.InsertAfter("TaskName")
[void]$tr1.InsertAfter(" VARCHAR(50) NOT NULL, ")
[void]$tr1.InsertAfter("StartDate")
[void]$tr1.InsertAfter(" DATE, ")
[void]$tr1.InsertAfter("DeadLine")
[void]$tr1.InsertAfter(" Date, Budget INT, Done BOOLEAN, FOREIGN KEY (PID) REFRENCES Project(PID) ON DELETE CASCADE);")
$tb1.TextFrame.WordWrap = -1
$tb1.TextFrame.AutoSize = 1

# --- TextBox 4: Project table -------------------------------------------
$tb2 = $s.Shapes.AddTextbox(1, 757282 / $EMU, 692696 / $EMU, 7704856 / $EMU, 923330 / $EMU)
$tb2.Name = "TextBox 4"
$tb2.Fill.Visible = 0
$tr2 = $tb2.TextFrame.TextRange
$tr2.Text = "CREATE TABLE Task(PID INTEGER NOT NULL PRIMARY KEY AUTOINCREMENT, "
[void]$tr2.InsertAfter("ProjectName")
[void]$tr2.InsertAfter(" VARCHAR(50) NOT NULL, ")
[void]$tr2.InsertAfter("StartDate")
[void]$tr2.InsertAfter(" DATE, ")
[void]$tr2.InsertAfter("DeadLine")
[void]$tr2.InsertAfter(" Date, Budget INT, ")
[void]$tr2.InsertAfter("Done BOOLEAN")
[void]$tr2.InsertAfter(");")
$tb2.TextFrame.WordWrap = -1
$tb2.TextFrame.AutoSize = 1

# --- TextBox 5: Resource table -------------------------------------------
$tb3 = $s.Shapes.AddTextbox(1, 757282 / $EMU, 3789040 / $EMU, 7704856 / $EMU, 923330 / $EMU)
$tb3.Name = "TextBox 5"
$tb3.Fill.Visible = 0
$tr3 = $tb3.TextFrame.TextRange
$tr3.Text = "CREATE TABLE Resource(RID INTEGER NOT NULL PRIMARY KEY AUTOINCREMENT, "
[void]$tr3.InsertAfter("FirstName")
[void]$tr3.InsertAfter(" ")
[void]$tr3.InsertAfter("VARCHAR(50), ")
[void]$tr3.InsertAfter("LastName")
[void]$tr3.InsertAfter(" VARCHAR(50), Salary INT,  ")
[void]$tr3.InsertAfter("UserName")
[void]$tr3.InsertAfter(" VARCHAR(10), Password char(32), Manager BOOLEAN);")
$tb3.TextFrame.WordWrap = -1
$tb3.TextFrame.AutoSize = 1

# --- TextBox 6: AssignedTask table ----------------------------------------
$tb4 = $s.Shapes.AddTextbox(1, 909682 / $EMU, 5013176 / $EMU, 7704856 / $EMU, 646331 / $EMU)
$tb4.Name = "TextBox 6"
$tb4.Fill.Visible = 0
$tr4 = $tb4.TextFrame.TextRange
$tr4.Text = "CREATE TABLE "
[void]$tr4.InsertAfter("AssignedTask")
[void]$tr4.InsertAfter("(RID INT NOT NULL, TID INT ")
[void]$tr4.InsertAfter("NOT NULL, ")
[void]$tr4.InsertAfter("Description TEXT, Deadline DATE, ")
[void]$tr4.InsertAfter("RequiredHoursWork")
[void]$tr4.InsertAfter(" INT, Done BOOLEAN);")
$tb4.TextFrame.WordWrap = -1
$tb4.TextFrame.AutoSize = 1
